$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '91.799.60'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.111.40'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.48'
$ws.Range('E5').Value = '  +2.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '621.51'
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.14'
$ws.Range('E7').Value = '  +5.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.372'
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('B10').Value = 'LidoStakedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.117.63'
$ws.Range('E10').Value = '  -0.53%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.762'
$ws.Range('E11').Value = '  +5.07%  '
$ws.Range('E12').Value = '  +3.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.51'
$ws.Range('E14').Value = '  -3.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.640.11'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.49'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.687.60'
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.112.53'
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.72'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.60'
$ws.Range('E20').Value = '  +2.20%  '
$ws.Range('B21').Value = 'PEPE'
$ws.Range('C21').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000215'
$ws.Range('E21').Value = '  +0.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.80'
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '448.92'
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.22'
$ws.Range('E24').Value = '  +2.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.88'
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '90.73'
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.02'
$ws.Range('E27').Value = '  -3.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.281.15'
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.186'
$ws.Range('E30').Value = '  +15.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.239'
$ws.Range('E31').Value = '  +20.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.43'
$ws.Range('E32').Value = '  -2.82%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.171'
$ws.Range('E33').Value = '  +13.95%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.113'
$ws.Range('E35').Value = '  +32.03%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.65'
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.69'
$ws.Range('E37').Value = '  +7.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.17'
$ws.Range('E38').Value = '  +24.59%  '
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '494.18'
$ws.Range('E40').Value = '  -3.75%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.62'
$ws.Range('E41').Value = '  -5.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.30'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.424'
$ws.Range('E43').Value = '  +1.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.17'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.701'
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '154.39'
$ws.Range('E48').Value = '  +2.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.57'
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.35'
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '44.63'
$ws.Range('E51').Value = '  -2.83%  '
